$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 160.41176
$ws.Range("I33").Value = 157.9375
$ws.Range("J33").Value = 200
$ws.Range("K33").Value = 157.9375
$ws.Range("L33").Value = 200
$ws.Range("M33").Value = 71.0625
$ws.Range("N33").Value = -658
$ws.Range("H51").Value = 2225
$ws.Range("I51").Value = 1450
$ws.Range("J51").Value = 3000
$ws.Range("K51").Value = 1450
$ws.Range("L51").Value = 3000
$ws.Range("M51").Value = -966
$ws.Range("N51").Value = -3968
$ws.Range("H64").Value = 5986.041
$ws.Range("J64").Value = 7691.857
$ws.Range("L64").Value = 7691.857
$ws.Range("N64").Value = -8187.857
$ws.Range("H67").Value = 5986.041
$ws.Range("J67").Value = 7691.857
$ws.Range("L67").Value = 7691.857
$ws.Range("N67").Value = -9407.857
$ws.Range("H80").Value = 8370.267
$ws.Range("I80").Value = 479.16666
$ws.Range("J80").Value = 13631
$ws.Range("K80").Value = 1437.49998
$ws.Range("L80").Value = 40893
$ws.Range("M80").Value = -439.4999800000001
$ws.Range("N80").Value = -42889
$ws.Range("H83").Value = 8370.267
$ws.Range("I83").Value = 479.16666
$ws.Range("J83").Value = 13631
$ws.Range("K83").Value = 4312.49994
$ws.Range("L83").Value = 122679
$ws.Range("M83").Value = 679.5000600000003
$ws.Range("N83").Value = -132663
$ws.Range("H92").Value = 564.1429000000001
$ws.Range("I92").Value = 612.5
$ws.Range("J92").Value = 499.66666
$ws.Range("K92").Value = 612.5
$ws.Range("L92").Value = 499.66666
$ws.Range("M92").Value = 635.5
$ws.Range("N92").Value = -2995.66666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7028.3887
$ws.Range("I32").Value = 2321.4285
$ws.Range("K32").Value = 2321.4285
$ws.Range("M32").Value = -2034.4285
$ws.Range("H61").Value = 47702.727
$ws.Range("I61").Value = 2208.5264
$ws.Range("K61").Value = 2208.5264
$ws.Range("M61").Value = -1996.5264
$ws.Range("H63").Value = 1863.5454
$ws.Range("I63").Value = 1924.9
$ws.Range("K63").Value = 1924.9
$ws.Range("M63").Value = -1238.9
$ws.Range("H66").Value = 1863.5454
$ws.Range("I66").Value = 1924.9
$ws.Range("K66").Value = 9624.5
$ws.Range("M66").Value = -6192.5
$ws.Range("H76").Value = 333403330
$ws.Range("J76").Value = 333403330
$ws.Range("L76").Value = 333403330
$ws.Range("N76").Value = -333404006
$ws.Range("H79").Value = 333403330
$ws.Range("J79").Value = 333403330
$ws.Range("L79").Value = 333403330
$ws.Range("N79").Value = -333405670
$ws.Range("H97").Value = 689.1053000000001
$ws.Range("I97").Value = 616.2778
$ws.Range("J97").Value = 2000
$ws.Range("K97").Value = 616.2778
$ws.Range("L97").Value = 2000
$ws.Range("M97").Value = -120.2778
$ws.Range("N97").Value = -2992
$ws.Range("H128").Value = 67900
$ws.Range("J128").Value = 67900
$ws.Range("L128").Value = 67900
$ws.Range("N128").Value = -77860
$ws.Range("H132").Value = 2231.6667
$ws.Range("I132").Value = 2160
$ws.Range("K132").Value = 6480
$ws.Range("M132").Value = -3950
$ws.Range("H136").Value = 47702.727
$ws.Range("I136").Value = 2208.5264
$ws.Range("K136").Value = 6625.5792
$ws.Range("M136").Value = -4075.5792

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2626.65
$ws.Range("I86").Value = 2692.7856
$ws.Range("J86").Value = 2472.3333
$ws.Range("K86").Value = 2692.7856
$ws.Range("L86").Value = 2472.3333
$ws.Range("M86").Value = -1569.7856
$ws.Range("N86").Value = -4718.3333
$ws.Range("H89").Value = 2626.65
$ws.Range("I89").Value = 2692.7856
$ws.Range("J89").Value = 2472.3333
$ws.Range("K89").Value = 13463.928
$ws.Range("L89").Value = 12361.6665
$ws.Range("M89").Value = -7847.928
$ws.Range("N89").Value = -23593.6665

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 12500251
$ws.Range("I25").Value = 268.33334
$ws.Range("K25").Value = 805.0000200000001
$ws.Range("M25").Value = -636.0000200000001
$ws.Range("H30").Value = 12500251
$ws.Range("I30").Value = 268.33334
$ws.Range("K30").Value = 805.0000200000001
$ws.Range("M30").Value = -703.0000200000001
$ws.Range("H40").Value = 57.6
$ws.Range("I40").Value = 60.333332
$ws.Range("J40").Value = 53.5
$ws.Range("K40").Value = 241.333328
$ws.Range("L40").Value = 214
$ws.Range("M40").Value = -172.333328
$ws.Range("N40").Value = -352
$ws.Range("H136").Value = 3768.5
$ws.Range("I136").Value = 1716
$ws.Range("K136").Value = 5148
$ws.Range("M136").Value = -48

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 6304703
$ws.Range("I11").Value = 2256569.2
$ws.Range("J11").Value = 13388937
$ws.Range("K11").Value = 2256569.2
$ws.Range("L11").Value = 13388937
$ws.Range("M11").Value = -2256430.2
$ws.Range("N11").Value = -13389215
$ws.Range("H80").Value = 4184.5
$ws.Range("I80").Value = 3251.75
$ws.Range("J80").Value = 6050
$ws.Range("K80").Value = 3251.75
$ws.Range("L80").Value = 6050
$ws.Range("M80").Value = -2253.75
$ws.Range("N80").Value = -8046
$ws.Range("H83").Value = 4184.5
$ws.Range("I83").Value = 3251.75
$ws.Range("J83").Value = 6050
$ws.Range("K83").Value = 16258.75
$ws.Range("L83").Value = 30250
$ws.Range("M83").Value = -11266.75
$ws.Range("N83").Value = -40234
$ws.Range("H131").Value = 92061
$ws.Range("J131").Value = 92061
$ws.Range("L131").Value = 92061
$ws.Range("N131").Value = -102141

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 15131.286
$ws.Range("I22").Value = 669
$ws.Range("J22").Value = 17541.666
$ws.Range("K22").Value = 669
$ws.Range("L22").Value = 17541.666
$ws.Range("M22").Value = -374
$ws.Range("N22").Value = -18131.666
$ws.Range("H27").Value = 15131.286
$ws.Range("I27").Value = 669
$ws.Range("J27").Value = 17541.666
$ws.Range("K27").Value = 669
$ws.Range("L27").Value = 17541.666
$ws.Range("M27").Value = -562
$ws.Range("N27").Value = -17755.666
$ws.Range("H53").Value = 7546.6665
$ws.Range("I53").Value = 4820
$ws.Range("J53").Value = 13000
$ws.Range("K53").Value = 4820
$ws.Range("L53").Value = 13000
$ws.Range("M53").Value = -4302
$ws.Range("N53").Value = -14036
$ws.Range("H82").Value = 1953.4
$ws.Range("I82").Value = 1234
$ws.Range("J82").Value = 2433
$ws.Range("K82").Value = 1234
$ws.Range("L82").Value = 2433
$ws.Range("M82").Value = -873
$ws.Range("N82").Value = -3155
$ws.Range("H85").Value = 1953.4
$ws.Range("I85").Value = 1234
$ws.Range("J85").Value = 2433
$ws.Range("K85").Value = 1234
$ws.Range("L85").Value = 2433
$ws.Range("M85").Value = 14
$ws.Range("N85").Value = -4929
$ws.Range("H93").Value = 2721.8333
$ws.Range("I93").Value = 2466.6
$ws.Range("K93").Value = 2466.6
$ws.Range("M93").Value = -1218.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H53").Value = 60084
$ws.Range("J53").Value = 60084
$ws.Range("L53").Value = 60084
$ws.Range("N53").Value = -61298
$ws.Range("H61").Value = 1847223.9
$ws.Range("I61").Value = 2026940.6
$ws.Range("K61").Value = 2026940.6
$ws.Range("M61").Value = -2026648.6
$ws.Range("H81").Value = 2189.9333
$ws.Range("I81").Value = 638.55554
$ws.Range("J81").Value = 4517
$ws.Range("K81").Value = 1277.11108
$ws.Range("L81").Value = 9034
$ws.Range("M81").Value = -216.1110799999999
$ws.Range("N81").Value = -11156
$ws.Range("H84").Value = 2189.9333
$ws.Range("I84").Value = 638.55554
$ws.Range("J84").Value = 4517
$ws.Range("K84").Value = 6385.555399999999
$ws.Range("L84").Value = 45170
$ws.Range("M84").Value = -1081.555399999999
$ws.Range("N84").Value = -55778
$ws.Range("H96").Value = 6606709.5
$ws.Range("I96").Value = 36183.168
$ws.Range("J96").Value = 26318288
$ws.Range("K96").Value = 36183.168
$ws.Range("L96").Value = 26318288
$ws.Range("M96").Value = -34810.168
$ws.Range("N96").Value = -26321034
$ws.Range("H100").Value = 4202569
$ws.Range("J100").Value = 783.25
$ws.Range("L100").Value = 1566.5
$ws.Range("N100").Value = -2648.5
$ws.Range("H101").Value = 46899
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 46899
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 46899
$ws.Range("M101").ClearContents()
$ws.Range("N101").Value = -53389
$ws.Range("H113").Value = 2252.0908
$ws.Range("I113").Value = 2816.1667
$ws.Range("K113").Value = 8448.500100000001
$ws.Range("M113").Value = -6278.500100000001
$ws.Range("H123").Value = 74790
$ws.Range("J123").Value = 74790
$ws.Range("L123").Value = 74790
$ws.Range("N123").Value = -84590
$ws.Range("H136").Value = 1351.6666
$ws.Range("I136").Value = 1159.2593
$ws.Range("K136").Value = 3477.7779
$ws.Range("M136").Value = -927.7779
